$wb = $excel.ActiveWorkbook

# --- Metadata sheet updates ---
$meta = $wb.Worksheets.Item("Metadata")

# Version: 5.0.0 -> 6.0.0
$meta.Range("B3").Value = "6.0.0"

# Date: updated publication timestamp
$meta.Range("B8").Value = "2022-01-21T20:46:54+00:00"

# Publisher value was previously blank, now set
$meta.Range("B9").Value = "Alvearie Team"

# The duplicated "Contact" row (row 10) becomes "Jurisdiction" / "United States of America"
$meta.Range("A10").Value = "Jurisdiction"
$meta.Range("B10").Value = "United States of America"

# The second duplicate "Contact" row (row 11) is removed entirely, shifting later rows up
$meta.Rows.Item(11).Delete()

# --- Elements sheet updates ---
$elements = $wb.Worksheets.Item("Elements")

# Root Extension row: Short/Definition now describe the Compound Code extension
$elements.Range("K2").Value = "Compound Code"
$elements.Range("L2").Value = "Customer-specific code for the compound of the drug"
